$wb = $excel.ActiveWorkbook

# --- Sheet "Obra": normalize dates from DD-MM-YYYY to YYYY-MM-DD ---
$wsObra = $wb.Worksheets.Item("Obra")
$wsObra.Range("C2").Value = "1965-01-01"
$wsObra.Range("C3").Value = "1997-01-01"
$wsObra.Range("C4").Value = "1983-01-01"

# --- Sheet "Referentes": normalize dates, strip stray quotes, fill blanks with NULL ---
$wsRef = $wb.Worksheets.Item("Referentes")
$wsRef.Range("B2").Value = "Doble suicidio en El Sisga"
$wsRef.Range("C2").Value = "1965-06-29"
$wsRef.Range("C3").Value = "1996-05-24"
$wsRef.Range("C4").Value = "NULL"
$wsRef.Range("D4").Value = "NULL"
$wsRef.Range("C5").Value = "NULL"
$wsRef.Range("D5").Value = "NULL"

# --- Selections / active sheet to match the author's final UI state ---
$wsObra.Range("C5").Select()
$wsRef.Activate()
$wsRef.Range("B5").Select()
